$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Feedback" header in G1, matching the style used by the other headers (F1)
$ws.Range("G1").Value = "Feedback"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Update the MarksAwarded values that changed alongside the new feedback
$ws.Range("F6").Value = 7.5
$ws.Range("F7").Value = 5.5
$ws.Range("F8").Value = 9
$ws.Range("F10").Value = 6.5

# Populate the Feedback column for the free-response questions (rows 2-10)
$ws.Range("G2").Value = 'None'
$ws.Range("G3").Value = 'Your answer is correct, but it lacks explanation and clarity. To improve, provide a brief explanation of how you derived the velocity equation (ds/dt = 6t + 2) and how you applied it to find the velocity at t = 2s. Additionally, consider adding units to your answer to make it more precise.'
$ws.Range("G4").Value = 'None'
$ws.Range("G5").Value = 'Your answer is correct, but it''s a straightforward copy of the correct answer. To improve, try to explain the reasoning behind the calculation. For example, you could mention that the initial velocity is used to calculate the maximum height using the equation u^2 / (2g), and then plug in the given values to get the answer. This will demonstrate your understanding of the concept and not just the calculation.'
$ws.Range("G6").Value = 'Your answer is concise and correctly states the equation of motion. However, it lacks the actual derivation process. To improve, provide the step-by-step derivation of the equation, starting from the initial velocity (u), acceleration (a), and final velocity (v). This will demonstrate your understanding of the concept and show how the equation is derived.'
$ws.Range("G7").Value = 'Your answer correctly defines average speed as the total distance traveled divided by the total time taken. However, it lacks information about instantaneous speed, which is a crucial concept in understanding speed. To improve, expand your answer to include the definition of instantaneous speed as the speed at any given instant of time. This will provide a more comprehensive understanding of the difference between the two concepts.'
$ws.Range("G8").Value = 'Your answer is close, but it''s missing the calculation steps. To improve, show your work and calculate the position using the equation s = ut + 0.5at^2. In this case, you would multiply the initial velocity (5 m/s) by the time (4 s) to get 20 m, then add half the acceleration (1.5 m/s^2) multiplied by the time squared (16) to get 24 m. Finally, add these two values to get the total position of 44 m.'
$ws.Range("G9").Value = 'Your answer is mostly correct, but it lacks the conversion of the average speed from km/hr to m/s. To improve, include the step where you convert the average speed from km/hr to m/s by multiplying by 1000/3600. This will ensure that your answer is accurate and in the correct units.'
$ws.Range("G10").Value = 'Your answer correctly identifies the reason for the same time of flight, citing the constant acceleration due to gravity. However, it lacks a clear explanation of why the acceleration due to gravity is independent of mass. To improve, expand on this point by stating that the acceleration due to gravity (g) is a fundamental constant, unaffected by the mass of the object. This highlights the underlying physics behind the phenomenon.'

# Rows 11-46 correspond to multiple-choice questions with no written feedback;
# touch their formatting (a no-op) so the Feedback cells still materialize as
# present-but-empty cells in that column, matching the source report.
$ws.Range("G11:G46").Font.Bold = $false

